$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '23.452.07'
$ws.Range('E2').Value = '  +1.16%  '
Set-TextValue $ws.Range('D3') '1.638.94'
$ws.Range('E3').Value = '  +2.33%  '
Set-TextValue $ws.Range('D4') '1.001'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +0.01%  '
Set-TextValue $ws.Range('D6') '306.15'
$ws.Range('E6').Value = '  +0.97%  '
Set-TextValue $ws.Range('D7') '0.3756'
$ws.Range('E7').Value = '  -0.58%  '
Set-TextValue $ws.Range('D8') '52.10'
$ws.Range('E8').Value = '  +0.13%  '
Set-TextValue $ws.Range('D9') '0.3640'
$ws.Range('E9').Value = '  +0.68%  '
Set-TextValue $ws.Range('D10') '1.262'
$ws.Range('E10').Value = '  -0.56%  '
Set-TextValue $ws.Range('D11') '0.08144'
$ws.Range('E11').Value = '  +0.35%  '
Set-TextValue $ws.Range('D12') '1.001'
$ws.Range('E12').Value = '  +0.05%  '
Set-TextValue $ws.Range('D13') '22.98'
$ws.Range('E13').Value = '  +1.05%  '
Set-TextValue $ws.Range('D14') '6.629'
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('E15').Value = '  +2.57%  '
Set-TextValue $ws.Range('D16') '7.368'
$ws.Range('E16').Value = '  -0.67%  '
Set-TextValue $ws.Range('D17') '1.635.71'
$ws.Range('E17').Value = '  +2.20%  '
Set-TextValue $ws.Range('D18') '94.57'
$ws.Range('E18').Value = '  +0.50%  '
Set-TextValue $ws.Range('D19') '0.06914'
$ws.Range('E19').Value = '  +0.61%  '
Set-TextValue $ws.Range('D20') '18.18'
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('E21').Value = '  -0.03%  '
Set-TextValue $ws.Range('D22') '1.0000'
$ws.Range('E22').Value = '  -0.05%  '
Set-TextValue $ws.Range('D23') '23.474.16'
$ws.Range('E23').Value = '  +1.22%  '
Set-TextValue $ws.Range('D24') '12.79'
$ws.Range('E24').Value = '  -1.52%  '
Set-TextValue $ws.Range('D25') '3.091'
$ws.Range('E25').Value = '  +3.77%  '
Set-TextValue $ws.Range('D26') '2.421'
$ws.Range('E26').Value = '  +0.83%  '
Set-TextValue $ws.Range('D27') '21.23'
$ws.Range('E27').Value = '  -0.05%  '
Set-TextValue $ws.Range('D28') '150.94'
$ws.Range('E28').Value = '  +1.07%  '
Set-TextValue $ws.Range('D29') '5.344'
$ws.Range('E29').Value = '  +1.72%  '
Set-TextValue $ws.Range('D30') '136.89'
$ws.Range('E30').Value = '  +2.24%  '
Set-TextValue $ws.Range('D31') '2.310'
$ws.Range('E31').Value = '  -2.90%  '
Set-TextValue $ws.Range('D32') '1.820.02'
$ws.Range('E32').Value = '  +2.32%  '
Set-TextValue $ws.Range('D33') '6.766'
$ws.Range('E33').Value = '  -0.09%  '
Set-TextValue $ws.Range('D34') '0.9649'
$ws.Range('E34').Value = '  -0.43%  '
Set-TextValue $ws.Range('D35') '0.02836'
$ws.Range('E35').Value = '  +4.34%  '
Set-TextValue $ws.Range('D36') '10.35'
$ws.Range('E36').Value = '  +0.48%  '
Set-TextValue $ws.Range('D37') '0.07306'
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('E38').Value = '  +1.05%  '
Set-TextValue $ws.Range('D39') '0.08839'
$ws.Range('E39').Value = '  +0.39%  '
Set-TextValue $ws.Range('D40') '6.126'
$ws.Range('E40').Value = '  +0.83%  '
Set-TextValue $ws.Range('D42') '0.7106'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D43') '16.27'
$ws.Range('E43').Value = '  +4.39%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D44') '12.47'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('E45').Value = '  +0.45%  '
Set-TextValue $ws.Range('D46') '2.340'
$ws.Range('E46').Value = '  +1.19%  '
Set-TextValue $ws.Range('D47') '0.9998'
$ws.Range('E47').Value = '  +0.00%  '
Set-TextValue $ws.Range('D48') '4.017'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E49').Value = '  +0.22%  '
Set-TextValue $ws.Range('D50') '128.81'
$ws.Range('E50').Value = '  -2.52%  '
Set-TextValue $ws.Range('D51') '1.207'
$ws.Range('E51').Value = '  +0.63%  '
